$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2620
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2620
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2620
$ws.Range("N51").Value = -3588

$ws.Range("H92").Value = 2052547.1
$ws.Range("I92").Value = 2462857.2
$ws.Range("J92").Value = 997
$ws.Range("K92").Value = 2462857.2
$ws.Range("L92").Value = 997
$ws.Range("M92").Value = -2461609.2
$ws.Range("N92").Value = -3493

$ws.Range("H132").Value = 928.4706
$ws.Range("I132").Value = 824.561
$ws.Range("J132").Value = 1354.5
$ws.Range("K132").Value = 2473.683
$ws.Range("L132").Value = 4063.5
$ws.Range("M132").Value = 56.31700000000001
$ws.Range("N132").Value = -9123.5

$ws.Range("H135").Value = 578.9286
$ws.Range("I135").Value = 508.75
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 4578.75
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -2043.75
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11628396
$ws.Range("I2").Value = 23255812
$ws.Range("J2").Value = 980
$ws.Range("K2").Value = 23255812
$ws.Range("L2").Value = 980
$ws.Range("M2").Value = -23255699
$ws.Range("N2").Value = -1206

$ws.Range("H32").Value = 3296.0881
$ws.Range("I32").Value = 1997.8596
$ws.Range("J32").Value = 10023.272
$ws.Range("K32").Value = 1997.8596
$ws.Range("L32").Value = 10023.272
$ws.Range("M32").Value = -1710.8596
$ws.Range("N32").Value = -10597.272

$ws.Range("H45").Value = 1692
$ws.Range("I45").Value = 1278
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 1278
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -901
$ws.Range("N45").Value = -2653

$ws.Range("H57").Value = 2000
$ws.Range("I57").Value = 2000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 2000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -1516

$ws.Range("H74").Value = 1676.8148
$ws.Range("I74").Value = 1237.0476
$ws.Range("J74").Value = 3216
$ws.Range("K74").Value = 1237.0476
$ws.Range("L74").Value = 3216
$ws.Range("M74").Value = -363.0476000000001
$ws.Range("N74").Value = -4964

$ws.Range("H77").Value = 1676.8148
$ws.Range("I77").Value = 1237.0476
$ws.Range("J77").Value = 3216
$ws.Range("K77").Value = 6185.238
$ws.Range("L77").Value = 16080
$ws.Range("M77").Value = -1817.238
$ws.Range("N77").Value = -24816

$ws.Range("H80").Value = 60000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 60000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -61996

$ws.Range("H83").Value = 60000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 60000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -189984

$ws.Range("H102").Value = 1485.6471
$ws.Range("I102").Value = 1217.0667
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 1217.0667
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = 404.9332999999999
$ws.Range("N102").Value = -6744

$ws.Range("H116").Value = 11628396
$ws.Range("I116").Value = 23255812
$ws.Range("J116").Value = 980
$ws.Range("K116").Value = 23255812
$ws.Range("L116").Value = 980
$ws.Range("M116").Value = -23253518
$ws.Range("N116").Value = -5568

$ws.Range("H122").Value = 1750
$ws.Range("I122").Value = 1750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5250
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2800

$ws.Range("H132").Value = 1826.1818
$ws.Range("I132").Value = 1231.9584
$ws.Range("J132").Value = 3410.7778
$ws.Range("K132").Value = 3695.8752
$ws.Range("L132").Value = 10232.3334
$ws.Range("M132").Value = -1165.8752
$ws.Range("N132").Value = -15292.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11628396
$ws.Range("I3").Value = 23255812
$ws.Range("J3").Value = 980
$ws.Range("K3").Value = 23255812
$ws.Range("L3").Value = 980
$ws.Range("M3").Value = -23255698
$ws.Range("N3").Value = -1208

$ws.Range("H99").Value = 1730
$ws.Range("I99").Value = 1386.6666
$ws.Range("J99").Value = 2171.4285
$ws.Range("K99").Value = 1386.6666
$ws.Range("L99").Value = 2171.4285
$ws.Range("M99").Value = 111.3334
$ws.Range("N99").Value = -5167.4285

$ws.Range("H105").Value = 2195.625
$ws.Range("I105").Value = 2187.0715
$ws.Range("J105").Value = 2255.5
$ws.Range("K105").Value = 2187.0715
$ws.Range("L105").Value = 2255.5
$ws.Range("M105").Value = -440.0715
$ws.Range("N105").Value = -5749.5

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 1217.8182
$ws.Range("I107").Value = 1089.6
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 1089.6
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = 830.4000000000001
$ws.Range("N107").Value = -6340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2845.3684
$ws.Range("I31").Value = 4270.6665
$ws.Range("J31").Value = 2578.125
$ws.Range("K31").Value = 4270.6665
$ws.Range("L31").Value = 2578.125
$ws.Range("M31").Value = -3975.6665
$ws.Range("N31").Value = -3168.125

$ws.Range("H34").Value = 2845.3684
$ws.Range("I34").Value = 4270.6665
$ws.Range("J34").Value = 2578.125
$ws.Range("K34").Value = 4270.6665
$ws.Range("L34").Value = 2578.125
$ws.Range("M34").Value = -4068.6665
$ws.Range("N34").Value = -2982.125

$ws.Range("H99").Value = 2521.125
$ws.Range("I99").Value = 2042.25
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2042.25
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -544.25
$ws.Range("N99").Value = -5996

$ws.Range("H105").Value = 2345.5715
$ws.Range("I105").Value = 2354.75
$ws.Range("J105").Value = 2333.3333
$ws.Range("K105").Value = 2354.75
$ws.Range("L105").Value = 2333.3333
$ws.Range("M105").Value = -607.75
$ws.Range("N105").Value = -5827.3333

$ws.Range("H126").Value = 2521.125
$ws.Range("I126").Value = 2042.25
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6126.75
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -3656.75
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 2131.5833
$ws.Range("I132").Value = 1277.381
$ws.Range("J132").Value = 3327.4666
$ws.Range("K132").Value = 3832.143
$ws.Range("L132").Value = 9982.399800000001
$ws.Range("M132").Value = -1302.143
$ws.Range("N132").Value = -15042.3998

$ws.Range("H134").Value = 1116.0541
$ws.Range("I134").Value = 1119.2778
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 3357.8334
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -822.8334000000004
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 8560.538
$ws.Range("I113").Value = 50500.5
$ws.Range("J113").Value = 935.0909
$ws.Range("K113").Value = 151501.5
$ws.Range("L113").Value = 2805.2727
$ws.Range("M113").Value = -149331.5
$ws.Range("N113").Value = -7145.2727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2084.2856
$ws.Range("I122").Value = 1547.75
$ws.Range("J122").Value = 2799.6667
$ws.Range("K122").Value = 4643.25
$ws.Range("L122").Value = 8399.000100000001
$ws.Range("M122").Value = -2193.25
$ws.Range("N122").Value = -13299.0001

$ws.Range("H132").Value = 2566036.2
$ws.Range("I132").Value = 3206754.8
$ws.Range("J132").Value = 3162
$ws.Range("K132").Value = 9620264.399999999
$ws.Range("L132").Value = 9486
$ws.Range("M132").Value = -9617734.399999999
$ws.Range("N132").Value = -14546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3507.647
$ws.Range("I7").Value = 2070
$ws.Range("J7").Value = 6143.3335
$ws.Range("K7").Value = 2070
$ws.Range("L7").Value = 6143.3335
$ws.Range("M7").Value = -1958
$ws.Range("N7").Value = -6367.3335

$ws.Range("H100").Value = 1814.2858
$ws.Range("I100").Value = 1616.6666
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1616.6666
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1075.6666
$ws.Range("N100").Value = -4082

$ws.Range("H122").Value = 12602.5
$ws.Range("I122").Value = 14834.667
$ws.Range("J122").Value = 10370.333
$ws.Range("K122").Value = 44504.001
$ws.Range("L122").Value = 31110.999
$ws.Range("M122").Value = -42054.001
$ws.Range("N122").Value = -36010.999

$ws.Range("H126").Value = 3507.647
$ws.Range("I126").Value = 2070
$ws.Range("J126").Value = 6143.3335
$ws.Range("K126").Value = 6210
$ws.Range("L126").Value = 18430.0005
$ws.Range("M126").Value = -3740
$ws.Range("N126").Value = -23370.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 260425.67
$ws.Range("I122").Value = 260425.67
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 781277.01
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -778827.01

$ws.Range("H132").Value = 1521.72
$ws.Range("I132").Value = 1159.7142
$ws.Range("J132").Value = 3422.25
$ws.Range("K132").Value = 3479.1426
$ws.Range("L132").Value = 10266.75
$ws.Range("M132").Value = -949.1425999999997
$ws.Range("N132").Value = -15326.75
